# Apply the UniformF-HW20 update:
#  - insert 4 new HKL-group columns worth of offset/grid data? (no - new ROW categories:
#    OffsetF, OffsetA, RD Single, TD Single) before the existing HexGrid row
#  - insert a new "1Pair-B" column among the HKL-group header columns
#  - extend the grid from 17 data columns (C:S) to 18 data columns (C:T)
#  - extend from 5 data rows (3:7) to 9 data rows (3:11)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: running index header, extend with one more column (18) ---
$ws.Cells.Item(1, 20).Value = 18
$ws.Range("S1").Copy()
$ws.Range("T1").PasteSpecial(-4122)

# --- Row 2: HKL / pair-group column headers (now 18 of them, C2:T2) ---
$headers = @("[1, 1, 0]","[2, 0, 0]","[2, 1, 1]","[2, 2, 0]","[3, 1, 0]","[2, 2, 2]","[3, 2, 1]","[4, 0, 0]","1Pair-A","1Pair-B","2Pairs-A","2Pairs-B","3Pairs-A","3Pairs-B","3Pairs-C","4Pairs","5A4F","MaxUnique")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(2, $i + 3).Value = $headers[$i]
}

# --- Extend existing data rows 3-7 with a value of 1 in the new column T ---
for ($r = 3; $r -le 7; $r++) {
    $ws.Cells.Item($r, 20).Value = 1
}

# Row 7 ("HexGrid-90degTilt5degRes") becomes "OffsetF" - the old category is
# re-inserted further down (new row 11) and a fresh "OffsetF" pattern takes
# its old slot.
$ws.Cells.Item(7, 2).Value = "OffsetF"

# --- New rows 8-11: additional scan-pattern categories ---
$newRows = @(
    @{ Index = 6; Name = "OffsetA" },
    @{ Index = 7; Name = "RD Single" },
    @{ Index = 8; Name = "TD Single" },
    @{ Index = 9; Name = "HexGrid-90degTilt5degRes" }
)

$targetRow = 8
foreach ($entry in $newRows) {
    $ws.Cells.Item($targetRow, 1).Value = $entry.Index
    $ws.Cells.Item($targetRow, 2).Value = $entry.Name
    for ($c = 3; $c -le 20; $c++) {
        $ws.Cells.Item($targetRow, $c).Value = 1
    }
    $targetRow++
}

# --- Copy the bold/centered/bordered format from column A (row 7) down to the new rows ---
$ws.Range("A7").Copy()
$ws.Range("A8:A11").PasteSpecial(-4122)
$excel.CutCopyMode = $false
